# Apply edits described by the commit "set game selection taskObject based on task"
# to the Blackbear-Consultants Deliverable_4 Sprint Backlog workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "D3-Sprint 5": fill in the previously-empty contribution-percentage
# row (row 6, columns G:M) with actual values, using Percent number formats.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("D3-Sprint 5")

$ws5.Range("G6").NumberFormat = "0%"
$ws5.Range("G6").Value = 0.075

$ws5.Range("H6").NumberFormat = "0%"
$ws5.Range("H6").Value = 0.15

$ws5.Range("I6").Value = 42.5

$ws5.Range("J6").NumberFormat = "0%"
$ws5.Range("J6").Value = 0.05

$ws5.Range("K6").NumberFormat = "0.00%"
$ws5.Range("K6").Value = 0.275

$ws5.Range("L6").NumberFormat = "0%"
$ws5.Range("L6").Value = 0

$ws5.Range("M6").NumberFormat = "0%"
$ws5.Range("M6").Value = 0.05

# ---------------------------------------------------------------------------
# Sheet "D4-Sprint 6": update "Story Points Completed" (D) and "Sprint #" (F)
# values for several backlog items.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("D4-Sprint 6")

$ws6.Range("D2").Value = 4
$ws6.Range("D3").Value = 2
$ws6.Range("D4").Value = 0.5
$ws6.Range("F4").Value = 4
$ws6.Range("D5").Value = 2
$ws6.Range("F5").Value = 4
$ws6.Range("D6").Value = 2
$ws6.Range("F6").Value = 4
$ws6.Range("D14").Value = 1
$ws6.Range("F14").Value = 4
$ws6.Range("D15").Value = 0.5
$ws6.Range("D16").Value = 2
$ws6.Range("F16").Value = 4
$ws6.Range("D20").Value = 2
$ws6.Range("F20").Value = 4
$ws6.Range("D21").Value = 2

# ---------------------------------------------------------------------------
# Update the active sheet / selection state: the workbook's active tab moves
# from "D4-Sprint 6" to "D3-Sprint 5".
# ---------------------------------------------------------------------------
[void]$ws6.Activate()
[void]$ws6.Range("L2").Select()

[void]$ws5.Activate()
[void]$ws5.Range("L7").Select()
